$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$ws.Cells.Item(2, 11).Value = "2025-10-28T03:01:35.696706+00:00"
$ws.Cells.Item(3, 11).Value = "2025-10-28T03:01:38.406631+00:00"
$ws.Cells.Item(4, 11).Value = "2025-10-28T03:01:38.406661+00:00"
$ws.Cells.Item(5, 11).Value = "2025-10-28T03:01:38.406680+00:00"
$ws.Cells.Item(6, 11).Value = "2025-10-28T03:01:38.406698+00:00"
$ws.Cells.Item(7, 11).Value = "2025-10-28T03:01:38.406714+00:00"
$ws.Cells.Item(8, 11).Value = "2025-10-28T03:01:40.783084+00:00"
$ws.Cells.Item(9, 11).Value = "2025-10-28T03:01:40.783117+00:00"
$ws.Cells.Item(10, 11).Value = "2025-10-28T03:01:43.497716+00:00"
$ws.Cells.Item(11, 11).Value = "2025-10-28T03:01:45.862860+00:00"
$ws.Cells.Item(12, 11).Value = "2025-10-28T03:01:45.862890+00:00"
$ws.Cells.Item(13, 11).Value = "2025-10-28T03:01:48.695365+00:00"
$ws.Cells.Item(14, 11).Value = "2025-10-28T03:01:48.695393+00:00"
$ws.Cells.Item(15, 11).Value = "2025-10-28T03:01:48.695411+00:00"
$ws.Cells.Item(16, 11).Value = "2025-10-28T03:01:48.695427+00:00"
$ws.Cells.Item(17, 11).Value = "2025-10-28T03:01:56.921274+00:00"
$ws.Cells.Item(18, 11).Value = "2025-10-28T03:01:59.722973+00:00"
$ws.Cells.Item(19, 11).Value = "2025-10-28T03:02:02.097136+00:00"
$ws.Cells.Item(20, 11).Value = "2025-10-28T03:02:04.388043+00:00"
$ws.Cells.Item(21, 11).Value = "2025-10-28T03:02:04.388071+00:00"
$ws.Cells.Item(22, 11).Value = "2025-10-28T03:02:07.165817+00:00"
$ws.Cells.Item(23, 11).Value = "2025-10-28T03:02:07.165847+00:00"
$ws.Cells.Item(24, 11).Value = "2025-10-28T03:02:07.165865+00:00"
$ws.Cells.Item(25, 11).Value = "2025-10-28T03:02:07.165881+00:00"
$ws.Cells.Item(26, 11).Value = "2025-10-28T03:02:07.165897+00:00"
$ws.Cells.Item(27, 11).Value = "2025-10-28T03:02:09.518310+00:00"
$ws.Cells.Item(28, 11).Value = "2025-10-28T03:02:14.603208+00:00"
$ws.Cells.Item(29, 11).Value = "2025-10-28T03:02:14.603239+00:00"
$ws.Cells.Item(30, 11).Value = "2025-10-28T03:02:14.603257+00:00"
$ws.Cells.Item(31, 11).Value = "2025-10-28T03:02:14.603274+00:00"
$ws.Cells.Item(32, 11).Value = "2025-10-28T03:02:17.388868+00:00"
$ws.Cells.Item(33, 11).Value = "2025-10-28T03:02:17.388899+00:00"
$ws.Cells.Item(34, 11).Value = "2025-10-28T03:02:19.724444+00:00"
$ws.Cells.Item(35, 11).Value = "2025-10-28T03:02:19.724472+00:00"
$ws.Cells.Item(36, 11).Value = "2025-10-28T03:02:19.724492+00:00"
$ws.Cells.Item(37, 11).Value = "2025-10-28T03:02:19.724508+00:00"
$ws.Cells.Item(38, 11).Value = "2025-10-28T03:02:19.724524+00:00"
$ws.Cells.Item(39, 11).Value = "2025-10-28T03:02:19.724544+00:00"
$ws.Cells.Item(40, 11).Value = "2025-10-28T03:02:19.724560+00:00"
$ws.Cells.Item(41, 11).Value = "2025-10-28T03:02:19.724575+00:00"
$ws.Cells.Item(42, 11).Value = "2025-10-28T03:02:19.724589+00:00"
$ws.Cells.Item(43, 11).Value = "2025-10-28T03:02:22.062228+00:00"
$ws.Cells.Item(44, 11).Value = "2025-10-28T03:02:22.062257+00:00"
$ws.Cells.Item(45, 11).Value = "2025-10-28T03:02:27.534599+00:00"
$ws.Cells.Item(46, 11).Value = "2025-10-28T03:02:30.242140+00:00"
$ws.Cells.Item(47, 11).Value = "2025-10-28T03:02:30.242169+00:00"
$ws.Cells.Item(48, 11).Value = "2025-10-28T03:02:30.242185+00:00"
$ws.Cells.Item(49, 11).Value = "2025-10-28T03:02:30.242204+00:00"
